$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.533.64'
$ws.Range('E2').Value = '  -0.84%  '

$ws.Range('D3').Value = '3.909.82'
$ws.Range('E3').Value = '  +4.02%  '

$ws.Range('E4').Value = '  +0.11%  '

$ws.Range('D5').Value = '''603.85'
$ws.Range('E5').Value = '  +0.20%  '

$ws.Range('D6').Value = '''165.18'
$ws.Range('E6').Value = '  -0.87%  '

$ws.Range('D7').Value = '3.907.28'
$ws.Range('E7').Value = '  +4.00%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').Value = '''0.527'
$ws.Range('E9').Value = '  -2.14%  '

$ws.Range('D10').Value = '''0.165'
$ws.Range('E10').Value = '  -3.53%  '

$ws.Range('D11').Value = '''6.40'
$ws.Range('E11').Value = '  +0.42%  '

$ws.Range('D12').Value = '''0.458'
$ws.Range('E12').Value = '  -0.22%  '

$ws.Range('D13').Value = '''36.84'
$ws.Range('E13').Value = '  -2.33%  '

$ws.Range('E14').Value = '  -0.93%  '

$ws.Range('D15').Value = '4.563.26'
$ws.Range('E15').Value = '  +4.03%  '

$ws.Range('D16').Value = '3.943.26'
$ws.Range('E16').Value = '  +4.78%  '

$ws.Range('D17').Value = '68.746.79'
$ws.Range('E17').Value = '  -0.53%  '

$ws.Range('D18').Value = '''7.40'

$ws.Range('E19').Value = '  -1.07%  '

$ws.Range('D20').Value = '''16.97'
$ws.Range('E20').Value = '  -4.15%  '

$ws.Range('D21').Value = '''11.17'
$ws.Range('E21').Value = '  -1.38%  '

$ws.Range('D22').Value = '''482.94'
$ws.Range('E22').Value = '  -1.66%  '

$ws.Range('B23').Value = 'Polygon'
$ws.Range('C23').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D23').Value = '''0.716'
$ws.Range('E23').Value = '  -1.48%  '

$ws.Range('B24').Value = 'PEPE'
$ws.Range('C24').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D24').Value = '''0.0000168'
$ws.Range('E24').Value = '  +12.70%  '

$ws.Range('D25').Value = '''84.36'
$ws.Range('E25').Value = '  -0.41%  '

$ws.Range('D26').Value = '''2.23'
$ws.Range('E26').Value = '  -1.68%  '

$ws.Range('D27').Value = '''12.00'
$ws.Range('E27').Value = '  -2.32%  '

$ws.Range('D28').Value = '''10.08'
$ws.Range('E28').Value = '  +0.08%  '

$ws.Range('E29').Value = '  -0.08%  '

$ws.Range('E30').Value = '  -1.08%  '

$ws.Range('D31').Value = '4.060.60'

$ws.Range('E32').Value = '  -3.15%  '

$ws.Range('E33').Value = '  -2.39%  '

$ws.Range('E34').Value = '  +0.59%  '

$ws.Range('D35').Value = '3.853.12'
$ws.Range('E35').Value = '  +3.90%  '

$ws.Range('E36').Value = '  -1.14%  '

$ws.Range('E37').Value = '  +2.79%  '

$ws.Range('D38').Value = '''0.139'
$ws.Range('E38').Value = '  -0.08%  '

$ws.Range('E39').Value = '  -1.28%  '

$ws.Range('E40').Value = '  +0.06%  '

$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '''0.318'
$ws.Range('E41').Value = '  -2.22%  '

$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').Value = '''3.00'
$ws.Range('E42').Value = '  -3.45%  '

$ws.Range('D43').Value = '''434.56'
$ws.Range('E43').Value = '  +1.43%  '

$ws.Range('D44').Value = '''48.48'
$ws.Range('E44').Value = '  -0.22%  '

$ws.Range('E45').Value = '  -0.81%  '

$ws.Range('D47').Value = '''8.42'
$ws.Range('E47').Value = '  -0.18%  '

$ws.Range('D48').Value = '2.826.21'
$ws.Range('E48').Value = '  +0.58%  '

$ws.Range('D49').Value = '''141.85'
$ws.Range('E49').Value = '  -0.58%  '

$ws.Range('D50').Value = '''25.86'
$ws.Range('E50').Value = '  +8.35%  '

$ws.Range('D51').Value = '''0.0352'
$ws.Range('E51').Value = '  -0.04%  '
